# Apply the "effort.xlsx" update:
#  - Row 15: B15 becomes 4 (single combined effort entry), C15 cleared (no more split A/B effort)
#  - Row 16: B16 becomes 2.75, C16 becomes 1.25
#  - New row 25: 07/07/2013 entry, 2.5h, new task description
#  - Selection moves to A25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# --- Row 15 ---------------------------------------------------------------
$ws.Range("B15").Value2 = 4
$ws.Range("C15").ClearContents()

# --- Row 16 ---------------------------------------------------------------
$ws.Range("B16").Value2 = 2.75
$ws.Range("C16").Value2 = 1.25

# --- New row 25 -------------------------------------------------------------
$ws.Range("A25").Value2 = 41462
$ws.Range("A25").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B25").Value2 = 2.5
$ws.Range("D25").Value = "SVN branch: gcc versus g++. Revision of Makefile, support of Linux and Windows, modularization"

# --- Update the active selection to the newly added row -------------------
$ws.Range("A25").Select() | Out-Null
